# Added the ability to change the allowed load.
#
# Rows 4 and 5 ("168 часов" / "79 часов") get a user-editable number of
# cars (column D) instead of the hard-coded "4"; the max-files total in
# column E follows (max files per car * number of cars). Column F ("Факт
# нагрузка в %") shows the resulting load percentage, rounded to a whole
# number for display.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F3:F5 is a merged cell; the merged block as a whole recalculates to the
# same overall load (96%) because the new total for E3:E5 (1740+2520+2520)
# equals the old total (1740+3360+1680). Unmerge so every underlying cell
# in the block can carry its own (identical) text value, then write it.
$ws.Range("F3:F5").UnMerge()

$ws.Cells.Item(3, 6).NumberFormat = "@"
$ws.Cells.Item(3, 6).Value = "96"
$ws.Cells.Item(3, 6).Style = "Normal"

$ws.Cells.Item(4, 6).NumberFormat = "@"
$ws.Cells.Item(4, 6).Value = "96"
$ws.Cells.Item(4, 6).Style = "Normal"

$ws.Cells.Item(5, 6).NumberFormat = "@"
$ws.Cells.Item(5, 6).Value = "96"
$ws.Cells.Item(5, 6).Style = "Normal"

# Row 6 ("180 часов праздники/вых") load % rounded from 35.6 -> 36.
$ws.Cells.Item(6, 6).NumberFormat = "@"
$ws.Cells.Item(6, 6).Value = "36"
$ws.Cells.Item(6, 6).Style = "Normal"

# Row 7 ("180 часов ночь") load % rounded from 48.8 -> 49.
$ws.Cells.Item(7, 6).NumberFormat = "@"
$ws.Cells.Item(7, 6).Value = "49"
$ws.Cells.Item(7, 6).Style = "Normal"

# Row 4 ("168 часов"): allowed load changed from 4 cars to 3 cars.
# Максимальное кол-во файлов в месяц (840) * кол-во машин (3) = 2520.
$ws.Cells.Item(4, 4).Value = 3
$ws.Cells.Item(4, 5).Value = 2520

# Row 5 ("79 часов"): allowed load changed from 4 cars to 6 cars.
# Максимальное кол-во файлов в месяц (420) * кол-во машин (6) = 2520.
$ws.Cells.Item(5, 4).Value = 6
$ws.Cells.Item(5, 5).Value = 2520
